# feat: add 2022-Q1 data
#
# - Insert a new sheet "2022-Q1" (holdings detail) right after "2021-Q4"
#   and before "总计".
# - Rebuild the "总计" (summary) sheet with a new top row for 2022-Q1,
#   shifting the previously existing rows down.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$oldTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Remove the old "总计" sheet and recreate the two sheets in the
#    right order so the new sheet lands between "2021-Q4" and "总计"
#    (mirrors sheetId 3 -> new sheet, sheetId 4 -> recreated "总计").
# ---------------------------------------------------------------------
$oldTotal.Delete()

$q1Sheet = $wb.Worksheets.Add($null, $q4)
$q1Sheet.Name = "2022-Q1"

$totalSheet = $wb.Worksheets.Add($null, $q1Sheet)
$totalSheet.Name = "总计"

# Match the page-margin convention already used by the sibling sheets
# (0.75in left/right, 1in top/bottom, 0.5in header/footer).
foreach ($sheet in @($q1Sheet, $totalSheet)) {
    $sheet.PageSetup.LeftMargin = 54
    $sheet.PageSetup.RightMargin = 54
    $sheet.PageSetup.TopMargin = 72
    $sheet.PageSetup.BottomMargin = 72
    $sheet.PageSetup.HeaderMargin = 36
    $sheet.PageSetup.FooterMargin = 36
}

# ---------------------------------------------------------------------
# 2. Populate "2022-Q1" sheet with the fund holdings detail.
#    Columns B (fund code) and D:G (text-formatted numeric-looking
#    figures) must stay text, so force Text number format before
#    assigning the value (otherwise leading zeros / trailing zeros are
#    lost to numeric coercion).
# ---------------------------------------------------------------------
$q1Sheet.Range("B1").Value = "基金代码"
$q1Sheet.Range("C1").Value = "基金名称"
$q1Sheet.Range("D1").Value = "基金规模"
$q1Sheet.Range("E1").Value = "股票总仓位"
$q1Sheet.Range("F1").Value = "仓位占比"
$q1Sheet.Range("G1").Value = "持有市值(亿元)"
$q1Sheet.Range("H1").Value = "仓位排名"

$q1Sheet.Range("A2").Value = 0
$q1Sheet.Range("B2:G2").NumberFormat = "@"
$q1Sheet.Range("B2").Value = "002252"
$q1Sheet.Range("C2").Value = "融通成长30灵活配置混合"
$q1Sheet.Range("D2").Value = "1.70"
$q1Sheet.Range("E2").Value = "78.26"
$q1Sheet.Range("F2").Value = "2.82"
$q1Sheet.Range("G2").Value = "0.0479"
$q1Sheet.Range("H2").Value = 6

$q1Sheet.Range("A3").Value = 1
$q1Sheet.Range("B3:G3").NumberFormat = "@"
$q1Sheet.Range("B3").Value = "009387"
$q1Sheet.Range("C3").Value = "嘉实稳福混合A"
$q1Sheet.Range("D3").Value = "0.08"
$q1Sheet.Range("E3").Value = "34.71"
$q1Sheet.Range("F3").Value = "1.74"
$q1Sheet.Range("G3").Value = "0.0014"
$q1Sheet.Range("H3").Value = 7

$q1Sheet.Range("A4").Value = 2
$q1Sheet.Range("B4:G4").NumberFormat = "@"
$q1Sheet.Range("B4").Value = "009388"
$q1Sheet.Range("C4").Value = "嘉实稳福混合C"
$q1Sheet.Range("D4").Value = "0.01"
$q1Sheet.Range("E4").Value = "34.71"
$q1Sheet.Range("F4").Value = "1.74"
$q1Sheet.Range("G4").Value = "0.0002"
$q1Sheet.Range("H4").Value = 7

# Re-use the existing header / first-column / data-row formatting
# (style index 2 for header & col A, default style for the rest,
# already used on the "2021-Q4" sheet) instead of creating new styles.
$q4.Range("B1:H1").Copy()
$q1Sheet.Range("B1:H1").PasteSpecial(-4122)

$q4.Range("A2").Copy()
$q1Sheet.Range("A2:A4").PasteSpecial(-4122)

$q4.Range("B2:H2").Copy()
$q1Sheet.Range("B2:H2").PasteSpecial(-4122)
$q1Sheet.Range("B3:H3").PasteSpecial(-4122)
$q1Sheet.Range("B4:H4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Populate "总计" sheet: 2022-Q1 on top, then the previously existing
#    2021-Q4 / 2021-Q1 rows shifted down by one.
# ---------------------------------------------------------------------
$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.05

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.12

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q1"
$totalSheet.Range("C4").Value = 4
$totalSheet.Range("D4").Value = 0.06

$q4.Range("B1:D1").Copy()
$totalSheet.Range("B1:D1").PasteSpecial(-4122)

$q4.Range("A2").Copy()
$totalSheet.Range("A2:A4").PasteSpecial(-4122)

Write-Output "done"
